{"js": "// Fix small spelling errors / typos, per the commit message.\n\n// 1) \"With other words\" -> \"In other words\"\n{\n  const results = context.document.body.search(\"With other words\", { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"In other words\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 2) Merge the \"ReLu activation\" paragraph with the \"validation accuracy\" paragraph\n//    that follows it, fixing \"which exception\" -> \"with exception\" and\n//    \"Here sigmoid activation is used.\" -> \"For the output layers, sigmoid activation is used.\"\n{\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n\n  let targetIndex = -1;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(\"which exception of the output layers\") !== -1) {\n      targetIndex = i;\n      break;\n    }\n  }\n\n  if (targetIndex !== -1) {\n    const reluParagraph = paragraphs.items[targetIndex];\n    const nextParagraph = paragraphs.items[targetIndex + 1];\n    nextParagraph.load(\"text\");\n    await context.sync();\n\n    const mergedText =\n      \"The layers in all models use ReLu activation, with exception of the output layers. \" +\n      \"For the output layers, sigmoid activation is used. \" +\n      nextParagraph.text;\n\n    reluParagraph.getRange().insertText(mergedText, Word.InsertLocation.replace);\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n\n// 3) \"running more other programs\" -> \"running additional other programs\"\n{\n  const results = context.document.body.search(\"running more other programs\", { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"running additional other programs\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Fix small spelling errors / typos, per the commit message.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $findText\n    $rng.Find.Replacement.Text = $replaceText\n    $rng.Find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$replaceText, 2) | Out-Null\n}\n\n# 1) \"With other words\" -> \"In other words\"\nReplace-Text \"With other words\" \"In other words\"\n\n# 2) Merge the \"ReLu activation\" paragraph with the paragraph that follows it\n#    (deleting the paragraph mark between them), while fixing\n#    \"which exception\" -> \"with exception\" and\n#    \"Here sigmoid activation is used.\" -> \"For the output layers, sigmoid activation is used.\"\nReplace-Text \"which exception of the output layers. Here sigmoid activation is used.^pFor the evaluation\" \"with exception of the output layers. For the output layers, sigmoid activation is used. For the evaluation\"\n\n# 3) \"running more other programs\" -> \"running additional other programs\"\nReplace-Text \"running more other programs\" \"running additional other programs\"\n"}
